$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert two new "Prrafodelista" paragraphs right before the paragraph
#    that starts with "Comparación últimas compras" (ilvl=2 / numId=1):
#       - "Productos más vendidos"   (ilvl=3 -> ListLevelNumber=4)
#       - "Por departamentos "       (ilvl=4 -> ListLevelNumber=5)
# ---------------------------------------------------------------------------
$target = $d.Content.Duplicate
$target.Find.Execute("Comparación últimas compras", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $target.Paragraphs(1)
$anchorIndex = $anchorPara.Index

# Insert a blank paragraph right before the anchor; the new blank paragraph
# takes over the original anchor index and "Comparación..." shifts down by one.
$null = $anchorPara.Range.InsertParagraphBefore()
$newPara1 = $d.Paragraphs($anchorIndex)
$newPara1.Range.ListFormat.ListLevelNumber = 4
$newPara1.Range.Text = "Productos más vendidos"

# Insert the second blank paragraph before the (now shifted) anchor paragraph
# (anchor is now at $anchorIndex + 1); same index-reuse logic applies.
$anchorPara2 = $d.Paragraphs($anchorIndex + 1)
$null = $anchorPara2.Range.InsertParagraphBefore()
$newPara2 = $d.Paragraphs($anchorIndex + 1)
$newPara2.Range.ListFormat.ListLevelNumber = 5
$newPara2.Range.Text = "Por departamentos "

# ---------------------------------------------------------------------------
# 2) Replace the numbered list formatting on the trailing empty paragraph
#    with a plain left indent (2880 twips = 144 pt), removing <w:numPr/>.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2508F569" w14:textId="77777777" w:rsidR="008634C4" w:rsidRPr="008634C4" w:rsidRDefault="008634C4" w:rsidP="008634C4"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="2880"/><w:rPr><w:lang w:val="es-US"/></w:rPr></w:pPr></w:p>'
$lastPara.Range.InsertXML($lastXml)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
